$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (recalculated figures) ---
$ws.Cells.Item(2, 7).Value = 0.0179024360383068
$ws.Cells.Item(3, 7).Value = 0.0179024360383068
$ws.Cells.Item(8, 6).Value = 0.03934
$ws.Cells.Item(8, 7).Value = 0.0584062782858298
$ws.Cells.Item(8, 12).Value = 0.03667
$ws.Cells.Item(9, 6).Value = 0.03934
$ws.Cells.Item(9, 7).Value = 0.0584062782858298
$ws.Cells.Item(9, 12).Value = 0.03667
$ws.Cells.Item(10, 7).Value = 0.0187573786648454
$ws.Cells.Item(11, 7).Value = 0.0187573786648454
$ws.Cells.Item(18, 7).Value = 0.129332681437586
$ws.Cells.Item(18, 12).Value = 0.04384
$ws.Cells.Item(19, 7).Value = 0.129332681437586
$ws.Cells.Item(19, 12).Value = 0.04384
$ws.Cells.Item(20, 7).Value = 0.0206549018992367
$ws.Cells.Item(21, 7).Value = 0.0206549018992367
$ws.Cells.Item(22, 7).Value = 0.0476660851063999
$ws.Cells.Item(23, 7).Value = 0.0476660851063999
$ws.Cells.Item(28, 7).Value = 0.129332681437586
$ws.Cells.Item(28, 12).Value = 0.04384
$ws.Cells.Item(29, 7).Value = 0.129332681437586
$ws.Cells.Item(29, 12).Value = 0.04384
$ws.Cells.Item(30, 7).Value = 0.0217977431936627
$ws.Cells.Item(31, 7).Value = 0.0217977431936627
$ws.Cells.Item(32, 7).Value = 0.0418630071222556
$ws.Cells.Item(33, 7).Value = 0.0418630071222556
$ws.Cells.Item(38, 7).Value = 0.172286260771876
$ws.Cells.Item(39, 7).Value = 0.172286260771876
$ws.Cells.Item(40, 7).Value = 0.0217145431812969
$ws.Cells.Item(41, 7).Value = 0.0217145431812969
$ws.Cells.Item(42, 7).Value = 0.0388682707692226
$ws.Cells.Item(43, 7).Value = 0.0388682707692226
$ws.Cells.Item(48, 6).Value = 0.04915
$ws.Cells.Item(48, 7).Value = 0.164927701182146
$ws.Cells.Item(48, 12).Value = 0.04584
$ws.Cells.Item(49, 6).Value = 0.04915
$ws.Cells.Item(49, 7).Value = 0.164927701182146
$ws.Cells.Item(49, 12).Value = 0.04584
$ws.Cells.Item(56, 7).Value = 0.0389934078892811
$ws.Cells.Item(57, 7).Value = 0.0389934078892811
$ws.Cells.Item(62, 7).Value = 0.170238474639236
$ws.Cells.Item(63, 7).Value = 0.170238474639236
$ws.Cells.Item(70, 7).Value = 0.0198842557181509
$ws.Cells.Item(71, 7).Value = 0.0198842557181509
$ws.Cells.Item(84, 7).Value = 0.018474882223899
$ws.Cells.Item(85, 7).Value = 0.018474882223899
$ws.Cells.Item(100, 7).Value = 0.017776375154006
$ws.Cells.Item(101, 7).Value = 0.017776375154006
$ws.Cells.Item(150, 7).Value = 0.0104903984352706
$ws.Cells.Item(150, 12).Value = 0.00107
$ws.Cells.Item(151, 7).Value = 0.0104903984352706
$ws.Cells.Item(151, 12).Value = 0.00107
$ws.Cells.Item(162, 7).Value = 2.64435031499221
$ws.Cells.Item(162, 8).Value = 6.5
$ws.Cells.Item(162, 9).Value = 5.71949
$ws.Cells.Item(170, 6).Value = 0.00286
$ws.Cells.Item(170, 7).Value = 0.0099634514774269
$ws.Cells.Item(170, 12).Value = 0.00195
$ws.Cells.Item(171, 6).Value = 0.00286
$ws.Cells.Item(171, 7).Value = 0.0099634514774269
$ws.Cells.Item(171, 12).Value = 0.00195
$ws.Cells.Item(182, 7).Value = 2.21758091806234
$ws.Cells.Item(182, 8).Value = 6.5
$ws.Cells.Item(182, 9).Value = 5.25394
$ws.Cells.Item(190, 7).Value = 0.0133142881812812
$ws.Cells.Item(190, 12).Value = 0.00223
$ws.Cells.Item(191, 7).Value = 0.0133142881812812
$ws.Cells.Item(191, 12).Value = 0.00223
$ws.Cells.Item(202, 7).Value = 1.77680334155714
$ws.Cells.Item(202, 8).Value = 5.84645811800856
$ws.Cells.Item(210, 7).Value = 0.0138993006966159
$ws.Cells.Item(210, 12).Value = 0.0035
$ws.Cells.Item(211, 7).Value = 0.0138993006966159
$ws.Cells.Item(211, 12).Value = 0.0035
$ws.Cells.Item(222, 7).Value = 1.41642055367772
$ws.Cells.Item(230, 7).Value = 0.0140759359379076
$ws.Cells.Item(230, 12).Value = 0.00357
$ws.Cells.Item(231, 7).Value = 0.0140759359379076
$ws.Cells.Item(231, 12).Value = 0.00357
$ws.Cells.Item(242, 7).Value = 0.906855702637608
$ws.Cells.Item(250, 6).Value = 0.00776
$ws.Cells.Item(250, 7).Value = 0.0139465029550058
$ws.Cells.Item(251, 6).Value = 0.00776
$ws.Cells.Item(251, 7).Value = 0.0139465029550058
$ws.Cells.Item(270, 6).Value = 0.00776
$ws.Cells.Item(270, 7).Value = 0.0134757003924076
$ws.Cells.Item(270, 12).Value = 0.00713
$ws.Cells.Item(271, 6).Value = 0.00776
$ws.Cells.Item(271, 7).Value = 0.0134757003924076
$ws.Cells.Item(271, 12).Value = 0.00713
$ws.Cells.Item(289, 6).Value = 0.00776
$ws.Cells.Item(289, 7).Value = 0.0138597178071621
$ws.Cells.Item(289, 12).Value = 0.00713
$ws.Cells.Item(290, 6).Value = 0.00776
$ws.Cells.Item(290, 7).Value = 0.0138597178071621
$ws.Cells.Item(290, 12).Value = 0.00713

# --- Append new rows 300-318 (2019-2023 period) ---
# Row 300
$ws.Cells.Item(300, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(300, 2).Value = "ASPM"
$ws.Cells.Item(300, 3).Value = "D"
$ws.Cells.Item(300, 4).Value = "2019 - 2023"
$ws.Cells.Item(300, 5).Value = "RepSite"
$ws.Cells.Item(300, 6).Value = 0.12
$ws.Cells.Item(300, 7).Value = 0.1164
$ws.Cells.Item(300, 8).Value = 0.135
$ws.Cells.Item(300, 9).Value = 0.135
$ws.Cells.Item(300, 10).Value = ""
$ws.Cells.Item(300, 11).Value = ""
$ws.Cells.Item(300, 12).Value = 0.12
$ws.Cells.Item(300, 13).Value = 0.1343
$ws.Cells.Item(300, 14).Value = 0.135
$ws.Cells.Item(300, 15).Value = 1790083.556
$ws.Cells.Item(300, 16).Value = 5500787.423
$ws.Cells.Item(300, 17).Value = "Horowhenua District"
$ws.Cells.Item(300, 18).Value = "Waiopehu"
$ws.Cells.Item(300, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(300, 20).Value = "Hoki_1a"
$ws.Cells.Item(300, 21).Value = ""

# Row 301
$ws.Cells.Item(301, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(301, 2).Value = "DRP (95th Percentile)"
$ws.Cells.Item(301, 3).Value = "D"
$ws.Cells.Item(301, 4).Value = "2019 - 2023"
$ws.Cells.Item(301, 5).Value = "RepSite"
$ws.Cells.Item(301, 6).Value = 0.022
$ws.Cells.Item(301, 7).Value = 0.02742
$ws.Cells.Item(301, 8).Value = 0.19
$ws.Cells.Item(301, 9).Value = 0.057
$ws.Cells.Item(301, 10).Value = ""
$ws.Cells.Item(301, 11).Value = ""
$ws.Cells.Item(301, 12).Value = 0.022
$ws.Cells.Item(301, 13).Value = 0.034
$ws.Cells.Item(301, 14).Value = 0.0415
$ws.Cells.Item(301, 15).Value = 1790083.556
$ws.Cells.Item(301, 16).Value = 5500787.423
$ws.Cells.Item(301, 17).Value = "Horowhenua District"
$ws.Cells.Item(301, 18).Value = "Waiopehu"
$ws.Cells.Item(301, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(301, 20).Value = "Hoki_1a"
$ws.Cells.Item(301, 21).Value = "mg/L"

# Row 302
$ws.Cells.Item(302, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(302, 2).Value = "DRP (Median)"
$ws.Cells.Item(302, 3).Value = "D"
$ws.Cells.Item(302, 4).Value = "2019 - 2023"
$ws.Cells.Item(302, 5).Value = "RepSite"
$ws.Cells.Item(302, 6).Value = 0.022
$ws.Cells.Item(302, 7).Value = 0.02742
$ws.Cells.Item(302, 8).Value = 0.19
$ws.Cells.Item(302, 9).Value = 0.057
$ws.Cells.Item(302, 10).Value = ""
$ws.Cells.Item(302, 11).Value = ""
$ws.Cells.Item(302, 12).Value = 0.022
$ws.Cells.Item(302, 13).Value = 0.034
$ws.Cells.Item(302, 14).Value = 0.0415
$ws.Cells.Item(302, 15).Value = 1790083.556
$ws.Cells.Item(302, 16).Value = 5500787.423
$ws.Cells.Item(302, 17).Value = "Horowhenua District"
$ws.Cells.Item(302, 18).Value = "Waiopehu"
$ws.Cells.Item(302, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(302, 20).Value = "Hoki_1a"
$ws.Cells.Item(302, 21).Value = "mg/L"

# Row 303
$ws.Cells.Item(303, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(303, 2).Value = "E coli (>260)"
$ws.Cells.Item(303, 3).Value = "D"
$ws.Cells.Item(303, 4).Value = "2019 - 2023"
$ws.Cells.Item(303, 5).Value = "RepSite"
$ws.Cells.Item(303, 6).Value = 232.5
$ws.Cells.Item(303, 7).Value = 456.54
$ws.Cells.Item(303, 8).Value = 3100
$ws.Cells.Item(303, 9).Value = 1741
$ws.Cells.Item(303, 10).Value = 20
$ws.Cells.Item(303, 11).Value = 44
$ws.Cells.Item(303, 12).Value = 320
$ws.Cells.Item(303, 13).Value = 714
$ws.Cells.Item(303, 14).Value = 1400
$ws.Cells.Item(303, 15).Value = 1790083.556
$ws.Cells.Item(303, 16).Value = 5500787.423
$ws.Cells.Item(303, 17).Value = "Horowhenua District"
$ws.Cells.Item(303, 18).Value = "Waiopehu"
$ws.Cells.Item(303, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(303, 20).Value = "Hoki_1a"
$ws.Cells.Item(303, 21).Value = "% exceedances over 260/100 mL"

# Row 304
$ws.Cells.Item(304, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(304, 2).Value = "E coli (>540)"
$ws.Cells.Item(304, 3).Value = "C"
$ws.Cells.Item(304, 4).Value = "2019 - 2023"
$ws.Cells.Item(304, 5).Value = "RepSite"
$ws.Cells.Item(304, 6).Value = 232.5
$ws.Cells.Item(304, 7).Value = 456.54
$ws.Cells.Item(304, 8).Value = 3100
$ws.Cells.Item(304, 9).Value = 1741
$ws.Cells.Item(304, 10).Value = 20
$ws.Cells.Item(304, 11).Value = 44
$ws.Cells.Item(304, 12).Value = 320
$ws.Cells.Item(304, 13).Value = 714
$ws.Cells.Item(304, 14).Value = 1400
$ws.Cells.Item(304, 15).Value = 1790083.556
$ws.Cells.Item(304, 16).Value = 5500787.423
$ws.Cells.Item(304, 17).Value = "Horowhenua District"
$ws.Cells.Item(304, 18).Value = "Waiopehu"
$ws.Cells.Item(304, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(304, 20).Value = "Hoki_1a"
$ws.Cells.Item(304, 21).Value = "% exceedances over 540/100 mL"

# Row 305
$ws.Cells.Item(305, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(305, 2).Value = "E coli (Median)"
$ws.Cells.Item(305, 3).Value = "D"
$ws.Cells.Item(305, 4).Value = "2019 - 2023"
$ws.Cells.Item(305, 5).Value = "RepSite"
$ws.Cells.Item(305, 6).Value = 232.5
$ws.Cells.Item(305, 7).Value = 456.54
$ws.Cells.Item(305, 8).Value = 3100
$ws.Cells.Item(305, 9).Value = 1741
$ws.Cells.Item(305, 10).Value = 20
$ws.Cells.Item(305, 11).Value = 44
$ws.Cells.Item(305, 12).Value = 320
$ws.Cells.Item(305, 13).Value = 714
$ws.Cells.Item(305, 14).Value = 1400
$ws.Cells.Item(305, 15).Value = 1790083.556
$ws.Cells.Item(305, 16).Value = 5500787.423
$ws.Cells.Item(305, 17).Value = "Horowhenua District"
$ws.Cells.Item(305, 18).Value = "Waiopehu"
$ws.Cells.Item(305, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(305, 20).Value = "Hoki_1a"
$ws.Cells.Item(305, 21).Value = "E. coli/100 mL"

# Row 306
$ws.Cells.Item(306, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(306, 2).Value = "E coli (95th Percentile)"
$ws.Cells.Item(306, 3).Value = "E"
$ws.Cells.Item(306, 4).Value = "2019 - 2023"
$ws.Cells.Item(306, 5).Value = "RepSite"
$ws.Cells.Item(306, 6).Value = 232.5
$ws.Cells.Item(306, 7).Value = 456.54
$ws.Cells.Item(306, 8).Value = 3100
$ws.Cells.Item(306, 9).Value = 1741
$ws.Cells.Item(306, 10).Value = 20
$ws.Cells.Item(306, 11).Value = 44
$ws.Cells.Item(306, 12).Value = 320
$ws.Cells.Item(306, 13).Value = 714
$ws.Cells.Item(306, 14).Value = 1400
$ws.Cells.Item(306, 15).Value = 1790083.556
$ws.Cells.Item(306, 16).Value = 5500787.423
$ws.Cells.Item(306, 17).Value = "Horowhenua District"
$ws.Cells.Item(306, 18).Value = "Waiopehu"
$ws.Cells.Item(306, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(306, 20).Value = "Hoki_1a"
$ws.Cells.Item(306, 21).Value = "E. coli/100 mL"

# Row 307
$ws.Cells.Item(307, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(307, 2).Value = "MCI"
$ws.Cells.Item(307, 3).Value = "D"
$ws.Cells.Item(307, 4).Value = "2019 - 2023"
$ws.Cells.Item(307, 5).Value = "RepSite"
$ws.Cells.Item(307, 6).Value = 70.67
$ws.Cells.Item(307, 7).Value = 66.614
$ws.Cells.Item(307, 8).Value = 74.8
$ws.Cells.Item(307, 9).Value = 74.8
$ws.Cells.Item(307, 10).Value = ""
$ws.Cells.Item(307, 11).Value = ""
$ws.Cells.Item(307, 12).Value = 70.67
$ws.Cells.Item(307, 13).Value = 74.17
$ws.Cells.Item(307, 14).Value = 74.8
$ws.Cells.Item(307, 15).Value = 1790083.556
$ws.Cells.Item(307, 16).Value = 5500787.423
$ws.Cells.Item(307, 17).Value = "Horowhenua District"
$ws.Cells.Item(307, 18).Value = "Waiopehu"
$ws.Cells.Item(307, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(307, 20).Value = "Hoki_1a"
$ws.Cells.Item(307, 21).Value = ""

# Row 308
$ws.Cells.Item(308, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(308, 2).Value = "Ammoniacal-N (95th Percentile)"
$ws.Cells.Item(308, 3).Value = "A"
$ws.Cells.Item(308, 4).Value = "2019 - 2023"
$ws.Cells.Item(308, 5).Value = "RepSite"
$ws.Cells.Item(308, 6).Value = 0.00804
$ws.Cells.Item(308, 7).Value = 0.0138662120892913
$ws.Cells.Item(308, 8).Value = 0.138461538461538
$ws.Cells.Item(308, 9).Value = 0.04271
$ws.Cells.Item(308, 10).Value = ""
$ws.Cells.Item(308, 11).Value = ""
$ws.Cells.Item(308, 12).Value = 0.00699
$ws.Cells.Item(308, 13).Value = 0.01649
$ws.Cells.Item(308, 14).Value = 0.03346
$ws.Cells.Item(308, 15).Value = 1790083.556
$ws.Cells.Item(308, 16).Value = 5500787.423
$ws.Cells.Item(308, 17).Value = "Horowhenua District"
$ws.Cells.Item(308, 18).Value = "Waiopehu"
$ws.Cells.Item(308, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(308, 20).Value = "Hoki_1a"
$ws.Cells.Item(308, 21).Value = "mg NH4-N/L"

# Row 309
$ws.Cells.Item(309, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(309, 2).Value = "Ammoniacal-N (Median)"
$ws.Cells.Item(309, 3).Value = "A"
$ws.Cells.Item(309, 4).Value = "2019 - 2023"
$ws.Cells.Item(309, 5).Value = "RepSite"
$ws.Cells.Item(309, 6).Value = 0.00804
$ws.Cells.Item(309, 7).Value = 0.0138662120892913
$ws.Cells.Item(309, 8).Value = 0.138461538461538
$ws.Cells.Item(309, 9).Value = 0.04271
$ws.Cells.Item(309, 10).Value = ""
$ws.Cells.Item(309, 11).Value = ""
$ws.Cells.Item(309, 12).Value = 0.00699
$ws.Cells.Item(309, 13).Value = 0.01649
$ws.Cells.Item(309, 14).Value = 0.03346
$ws.Cells.Item(309, 15).Value = 1790083.556
$ws.Cells.Item(309, 16).Value = 5500787.423
$ws.Cells.Item(309, 17).Value = "Horowhenua District"
$ws.Cells.Item(309, 18).Value = "Waiopehu"
$ws.Cells.Item(309, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(309, 20).Value = "Hoki_1a"
$ws.Cells.Item(309, 21).Value = "mg NH4-N/L"

# Row 310
$ws.Cells.Item(310, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(310, 2).Value = "Nitrate-N (95th Percentile)"
$ws.Cells.Item(310, 3).Value = "D"
$ws.Cells.Item(310, 4).Value = "2019 - 2023"
$ws.Cells.Item(310, 5).Value = "RepSite"
$ws.Cells.Item(310, 6).Value = 9.94
$ws.Cells.Item(310, 7).Value = 9.5918
$ws.Cells.Item(310, 8).Value = 14.7
$ws.Cells.Item(310, 9).Value = 12.3
$ws.Cells.Item(310, 10).Value = ""
$ws.Cells.Item(310, 11).Value = ""
$ws.Cells.Item(310, 12).Value = 10.3
$ws.Cells.Item(310, 13).Value = 11.7
$ws.Cells.Item(310, 14).Value = 12.1
$ws.Cells.Item(310, 15).Value = 1790083.556
$ws.Cells.Item(310, 16).Value = 5500787.423
$ws.Cells.Item(310, 17).Value = "Horowhenua District"
$ws.Cells.Item(310, 18).Value = "Waiopehu"
$ws.Cells.Item(310, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(310, 20).Value = "Hoki_1a"
$ws.Cells.Item(310, 21).Value = "mg NO3-N/L"

# Row 311
$ws.Cells.Item(311, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(311, 2).Value = "Nitrate-N (Median)"
$ws.Cells.Item(311, 3).Value = "D"
$ws.Cells.Item(311, 4).Value = "2019 - 2023"
$ws.Cells.Item(311, 5).Value = "RepSite"
$ws.Cells.Item(311, 6).Value = 9.94
$ws.Cells.Item(311, 7).Value = 9.5918
$ws.Cells.Item(311, 8).Value = 14.7
$ws.Cells.Item(311, 9).Value = 12.3
$ws.Cells.Item(311, 10).Value = ""
$ws.Cells.Item(311, 11).Value = ""
$ws.Cells.Item(311, 12).Value = 10.3
$ws.Cells.Item(311, 13).Value = 11.7
$ws.Cells.Item(311, 14).Value = 12.1
$ws.Cells.Item(311, 15).Value = 1790083.556
$ws.Cells.Item(311, 16).Value = 5500787.423
$ws.Cells.Item(311, 17).Value = "Horowhenua District"
$ws.Cells.Item(311, 18).Value = "Waiopehu"
$ws.Cells.Item(311, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(311, 20).Value = "Hoki_1a"
$ws.Cells.Item(311, 21).Value = "mg NO3-N/L"

# Row 312
$ws.Cells.Item(312, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(312, 2).Value = "QMCI"
$ws.Cells.Item(312, 3).Value = "D"
$ws.Cells.Item(312, 4).Value = "2019 - 2023"
$ws.Cells.Item(312, 5).Value = "RepSite"
$ws.Cells.Item(312, 6).Value = 3.527
$ws.Cells.Item(312, 7).Value = 3.5694
$ws.Cells.Item(312, 8).Value = 4.35
$ws.Cells.Item(312, 9).Value = 4.35
$ws.Cells.Item(312, 10).Value = ""
$ws.Cells.Item(312, 11).Value = ""
$ws.Cells.Item(312, 12).Value = 3.527
$ws.Cells.Item(312, 13).Value = 4.091
$ws.Cells.Item(312, 14).Value = 4.35
$ws.Cells.Item(312, 15).Value = 1790083.556
$ws.Cells.Item(312, 16).Value = 5500787.423
$ws.Cells.Item(312, 17).Value = "Horowhenua District"
$ws.Cells.Item(312, 18).Value = "Waiopehu"
$ws.Cells.Item(312, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(312, 20).Value = "Hoki_1a"
$ws.Cells.Item(312, 21).Value = ""

# Row 313
$ws.Cells.Item(313, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(313, 2).Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Cells.Item(313, 3).Value = ""
$ws.Cells.Item(313, 4).Value = "2019 - 2023"
$ws.Cells.Item(313, 5).Value = "RepSite"
$ws.Cells.Item(313, 6).Value = 10.003
$ws.Cells.Item(313, 7).Value = 9.68572
$ws.Cells.Item(313, 8).Value = 14.708
$ws.Cells.Item(313, 9).Value = 12.308
$ws.Cells.Item(313, 10).Value = ""
$ws.Cells.Item(313, 11).Value = ""
$ws.Cells.Item(313, 12).Value = 10.336
$ws.Cells.Item(313, 13).Value = 11.805
$ws.Cells.Item(313, 14).Value = 12.07
$ws.Cells.Item(313, 15).Value = 1790083.556
$ws.Cells.Item(313, 16).Value = 5500787.423
$ws.Cells.Item(313, 17).Value = "Horowhenua District"
$ws.Cells.Item(313, 18).Value = "Waiopehu"
$ws.Cells.Item(313, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(313, 20).Value = "Hoki_1a"
$ws.Cells.Item(313, 21).Value = "g/m3"

# Row 314
$ws.Cells.Item(314, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(314, 2).Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Cells.Item(314, 3).Value = ""
$ws.Cells.Item(314, 4).Value = "2019 - 2023"
$ws.Cells.Item(314, 5).Value = "RepSite"
$ws.Cells.Item(314, 6).Value = 10.003
$ws.Cells.Item(314, 7).Value = 9.68572
$ws.Cells.Item(314, 8).Value = 14.708
$ws.Cells.Item(314, 9).Value = 12.308
$ws.Cells.Item(314, 10).Value = ""
$ws.Cells.Item(314, 11).Value = ""
$ws.Cells.Item(314, 12).Value = 10.336
$ws.Cells.Item(314, 13).Value = 11.805
$ws.Cells.Item(314, 14).Value = 12.07
$ws.Cells.Item(314, 15).Value = 1790083.556
$ws.Cells.Item(314, 16).Value = 5500787.423
$ws.Cells.Item(314, 17).Value = "Horowhenua District"
$ws.Cells.Item(314, 18).Value = "Waiopehu"
$ws.Cells.Item(314, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(314, 20).Value = "Hoki_1a"
$ws.Cells.Item(314, 21).Value = "g/m3"

# Row 315
$ws.Cells.Item(315, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(315, 2).Value = "Total Nitrogen (95th Percentile)"
$ws.Cells.Item(315, 3).Value = ""
$ws.Cells.Item(315, 4).Value = "2019 - 2023"
$ws.Cells.Item(315, 5).Value = "RepSite"
$ws.Cells.Item(315, 6).Value = 9.95
$ws.Cells.Item(315, 7).Value = 9.8044
$ws.Cells.Item(315, 8).Value = 13.8
$ws.Cells.Item(315, 9).Value = 12.5
$ws.Cells.Item(315, 10).Value = ""
$ws.Cells.Item(315, 11).Value = ""
$ws.Cells.Item(315, 12).Value = 10.1
$ws.Cells.Item(315, 13).Value = 11.6
$ws.Cells.Item(315, 14).Value = 12
$ws.Cells.Item(315, 15).Value = 1790083.556
$ws.Cells.Item(315, 16).Value = 5500787.423
$ws.Cells.Item(315, 17).Value = "Horowhenua District"
$ws.Cells.Item(315, 18).Value = "Waiopehu"
$ws.Cells.Item(315, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(315, 20).Value = "Hoki_1a"
$ws.Cells.Item(315, 21).Value = "g/m3"

# Row 316
$ws.Cells.Item(316, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(316, 2).Value = "Total Nitrogen (Median)"
$ws.Cells.Item(316, 3).Value = ""
$ws.Cells.Item(316, 4).Value = "2019 - 2023"
$ws.Cells.Item(316, 5).Value = "RepSite"
$ws.Cells.Item(316, 6).Value = 9.95
$ws.Cells.Item(316, 7).Value = 9.8044
$ws.Cells.Item(316, 8).Value = 13.8
$ws.Cells.Item(316, 9).Value = 12.5
$ws.Cells.Item(316, 10).Value = ""
$ws.Cells.Item(316, 11).Value = ""
$ws.Cells.Item(316, 12).Value = 10.1
$ws.Cells.Item(316, 13).Value = 11.6
$ws.Cells.Item(316, 14).Value = 12
$ws.Cells.Item(316, 15).Value = 1790083.556
$ws.Cells.Item(316, 16).Value = 5500787.423
$ws.Cells.Item(316, 17).Value = "Horowhenua District"
$ws.Cells.Item(316, 18).Value = "Waiopehu"
$ws.Cells.Item(316, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(316, 20).Value = "Hoki_1a"
$ws.Cells.Item(316, 21).Value = "g/m3"

# Row 317
$ws.Cells.Item(317, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(317, 2).Value = "Total Phosphorus (95th Percentile)"
$ws.Cells.Item(317, 3).Value = ""
$ws.Cells.Item(317, 4).Value = "2019 - 2023"
$ws.Cells.Item(317, 5).Value = "RepSite"
$ws.Cells.Item(317, 6).Value = 0.0515
$ws.Cells.Item(317, 7).Value = 0.08106
$ws.Cells.Item(317, 8).Value = 0.598
$ws.Cells.Item(317, 9).Value = 0.172
$ws.Cells.Item(317, 10).Value = ""
$ws.Cells.Item(317, 11).Value = ""
$ws.Cells.Item(317, 12).Value = 0.04
$ws.Cells.Item(317, 13).Value = 0.11
$ws.Cells.Item(317, 14).Value = 0.147
$ws.Cells.Item(317, 15).Value = 1790083.556
$ws.Cells.Item(317, 16).Value = 5500787.423
$ws.Cells.Item(317, 17).Value = "Horowhenua District"
$ws.Cells.Item(317, 18).Value = "Waiopehu"
$ws.Cells.Item(317, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(317, 20).Value = "Hoki_1a"
$ws.Cells.Item(317, 21).Value = "g/m3"

# Row 318
$ws.Cells.Item(318, 1).Value = "Arawhata Drain at Hokio Beach Road"
$ws.Cells.Item(318, 2).Value = "Total Phosphorus (Median)"
$ws.Cells.Item(318, 3).Value = ""
$ws.Cells.Item(318, 4).Value = "2019 - 2023"
$ws.Cells.Item(318, 5).Value = "RepSite"
$ws.Cells.Item(318, 6).Value = 0.0515
$ws.Cells.Item(318, 7).Value = 0.08106
$ws.Cells.Item(318, 8).Value = 0.598
$ws.Cells.Item(318, 9).Value = 0.172
$ws.Cells.Item(318, 10).Value = ""
$ws.Cells.Item(318, 11).Value = ""
$ws.Cells.Item(318, 12).Value = 0.04
$ws.Cells.Item(318, 13).Value = 0.11
$ws.Cells.Item(318, 14).Value = 0.147
$ws.Cells.Item(318, 15).Value = 1790083.556
$ws.Cells.Item(318, 16).Value = 5500787.423
$ws.Cells.Item(318, 17).Value = "Horowhenua District"
$ws.Cells.Item(318, 18).Value = "Waiopehu"
$ws.Cells.Item(318, 19).Value = "Lake Horowhenua"
$ws.Cells.Item(318, 20).Value = "Hoki_1a"
$ws.Cells.Item(318, 21).Value = "g/m3"

